$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.633.15'
$ws.Range("E2").Value = '  +5.57%  '
$ws.Range("D3").Value = '3.799.08'
$ws.Range("E3").Value = '  +23.11%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = '3.790.81'
$ws.Range("E7").Value = '  +22.83%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.77%  '
$ws.Range("E10").Value = '  +13.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.40'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.507'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.74'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.65%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000262'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.31%  '
$ws.Range("D15").Value = '4.435.78'
$ws.Range("E15").Value = '  +23.22%  '
$ws.Range("D16").Value = '3.799.65'
$ws.Range("E16").Value = '  +23.21%  '
$ws.Range("D17").Value = '70.796.70'
$ws.Range("E17").Value = '  +5.88%  '
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '525.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +23.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.749'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.90%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000124'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +33.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.92'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +13.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +16.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.22%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.19'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.344'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.19'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.08%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '3.159.72'
$ws.Range("E42").Value = '  +12.94%  '
$ws.Range("B43").Value = 'Cosmos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +8.33%  '
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '427.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.71%  '
$ws.Range("B45").Value = 'Arweave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0370'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.51%  '
$ws.Range("E51").Value = '  +0.01%  '
